# Updated cryptos list on Thu Aug 15 03:49:01 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    # Forces the cell to hold a text value even when the string looks like a
    # number (e.g. "521.87"), matching the original inlineStr content, then
    # restores the cell's original (default) style so no visible formatting
    # change is introduced.
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "58.393.81"
$ws.Range("E2").Value = "  -4.53%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.647.77"
$ws.Range("E3").Value = "  -1.80%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.20%  "

# Row 5 - BNB
Set-TextValue "D5" "521.87"
$ws.Range("E5").Value = "  -0.88%  "

# Row 6 - Solana
Set-TextValue "D6" "144.57"
$ws.Range("E6").Value = "  -0.94%  "

# Row 7 - USDC
Set-TextValue "D7" "1.00"
$ws.Range("E7").Value = "  +0.27%  "

# Row 8 - XRP
$ws.Range("E8").Value = "  -1.27%  "

# Row 9 - Toncoin
$ws.Range("E9").Value = "  +3.24%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  -3.42%  "

# Row 11 - Cardano
$ws.Range("E11").Value = "  -0.46%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +1.69%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "3.113.79"
$ws.Range("E13").Value = "  -1.73%  "

# Row 14 - WrappedBTC
Set-TextValue "D14" "58.397.91"
$ws.Range("E14").Value = "  -4.32%  "

# Row 15 - Avalanche
Set-TextValue "D15" "20.95"
$ws.Range("E15").Value = "  -2.21%  "

# Row 16 - ShibaInu
Set-TextValue "D16" "0.0000136"
$ws.Range("E16").Value = "  -1.78%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.658.25"
$ws.Range("E17").Value = "  -9.98%  "

# Row 18 - BitcoinCash
Set-TextValue "D18" "338.45"
$ws.Range("E18").Value = "  -3.31%  "

# Row 19 - Polkadot
Set-TextValue "D19" "4.39"
$ws.Range("E19").Value = "  -3.07%  "

# Row 20 - Chainlink
Set-TextValue "D20" "10.45"
$ws.Range("E20").Value = "  -1.56%  "

# Row 21 - Uniswap
Set-TextValue "D21" "6.31"
$ws.Range("E21").Value = "  -0.52%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  +0.33%  "

# Row 23 - Litecoin
Set-TextValue "D23" "64.39"
$ws.Range("E23").Value = "  +0.79%  "

# Row 24 - Polygon
Set-TextValue "D24" "0.425"
$ws.Range("E24").Value = "  +0.54%  "

# Row 25 - Kaspa
$ws.Range("E25").Value = "  -2.25%  "

# Row 26 - Binance-PegBSC-USD
$ws.Range("E26").Value = "  +0.58%  "

# Row 27 - PEPE
$ws.Range("D27").Value = "0.0₃0798"
$ws.Range("E27").Value = "  -2.99%  "

# Row 28 - InternetComputer(DFINITY)
Set-TextValue "D28" "7.14"
$ws.Range("E28").Value = "  -3.23%  "

# Row 29 - Aptos
$ws.Range("E29").Value = "  -3.50%  "

# Row 30 - USDe
$ws.Range("E30").Value = "  +0.12%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -0.86%  "

# Row 32 - Monero
Set-TextValue "D32" "152.60"
$ws.Range("E32").Value = "  +1.56%  "

# Row 33 - EthereumClassic
Set-TextValue "D33" "18.87"
$ws.Range("E33").Value = "  -1.99%  "

# Row 34 - NEARProtocol
$ws.Range("E34").Value = "  -3.19%  "

# Row 35 - ImmutableX
Set-TextValue "D35" "1.18"
$ws.Range("E35").Value = "  -5.38%  "

# Row 36 - SuiNetwork
Set-TextValue "D36" "0.910"
$ws.Range("E36").Value = "  -4.48%  "

# Row 37 - Fetch.AI
Set-TextValue "D37" "0.865"
$ws.Range("E37").Value = "  -2.25%  "

# Row 38 - OKB
Set-TextValue "D38" "36.81"
$ws.Range("E38").Value = "  -0.51%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  -5.28%  "

# Row 40 - Filecoin
$ws.Range("E40").Value = "  -1.21%  "

# Row 41 - FirstDigitalUSD
$ws.Range("E41").Value = "  +0.44%  "

# Row 42 - Mantle
Set-TextValue "D42" "0.608"
$ws.Range("E42").Value = "  -0.89%  "

# Row 43 - Bittensor
Set-TextValue "D43" "273.55"
$ws.Range("E43").Value = "  -3.97%  "

# Row 44 - Stellar
$ws.Range("E44").Value = "  -2.08%  "

# Row 45 - EnergySwap
Set-TextValue "D45" "19.49"
$ws.Range("E45").Value = "  -3.23%  "

# Rows 46 and 47 swap content (Hedera <-> WhiteBITCoin) with updated values
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D46" "10.63"
$ws.Range("E46").Value = "  +1.51%  "

$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue "D47" "0.0535"
$ws.Range("E47").Value = "  -1.21%  "

# Row 48 - Maker
Set-TextValue "D48" "2.044.84"
$ws.Range("E48").Value = "  -4.90%  "

# Row 49 - RenderToken
Set-TextValue "D49" "4.70"
$ws.Range("E49").Value = "  -4.61%  "

# Row 50 - VeChain
Set-TextValue "D50" "0.0228"
$ws.Range("E50").Value = "  -3.17%  "

# Row 51 - InjectiveProtocol
Set-TextValue "D51" "18.35"
$ws.Range("E51").Value = "  -4.05%  "
